$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Cells.Item(2, 1).Value = "281474991265517-1753110601926"
$ws.Cells.Item(2, 2).Value = "No Seat Belt"
$ws.Cells.Item(2, 3).Value = "2025-07-21T09:10:01.926"

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "281474991265517"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "144"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "52215890"
$ws.Cells.Item(2, 6).Style = "Normal"

$ws.Cells.Item(2, 7).Value = "MARIO BALLINAS"
$ws.Cells.Item(2, 8).Value = 20.68147431
$ws.Cells.Item(2, 9).Value = -103.46047671
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991265517/1753110599426/g9yG40HADs-camera-video-segment-driver-1753110601926.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFFWEYYQK%2F20250722%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250722T160114Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjENX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFA8Y0PgXhDtNu2%2BQHuQpEZwh0JskrXSr5GnqoS%2FtvyLAiEA14WPHxlFKgBdu70KIBqHbLoKFzI6fR80kkWgcedYHasq5gMI7v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDDKCurDpNn5HydYpgCq6A7tVZ4RPeUwBMWtl7mltXrtDVDjafXNo41hLbZ8mNNOiXR0bnL%2FSwLen7qeGKU1jmlnfe0DcMSjLMqX%2FXjtuB%2FBIKxNbO5YvcixVrqwuMVhKdK5s0KbZrzZlSEVpt5ir%2Bb2slBXuGpBE%2Fk1negUCLUvngM1WwVQ9lSLEU5eu75OS6bSJ6ZT%2BmYw9UmY%2BFof9%2FmWjTehN6qLDSMCeo4gcUg2nrH%2Bh620YdZ1sJcTegwmxxf5GU0ZUWuuQRIMRphd5xiv1K%2FJgPiDkAJCAqgBBIs2E5%2FfH6Fu8TXYTl9P%2FI0PT%2FZvy%2BN1i0HMilb6vKDTQsvZ4TcMzmSyVtvSwSaQzzZlpI3KgaEmhUiJKowNpXdgu88rXpBTMhyprbQH8P1MBr2n9ONmLbu0g4RjQZMRfUfey6uAfPbEur77F%2BTjhBEDM3mwejefwA3AeKZv7TiEgvbDOD9QYMxZV8XUTR9kMNlJCsrVR7pDR2CMC9%2F32ALrQIbas6UXYrRPgglwN9uobg59jHNiY72uGCJCe%2Bp0NIF6sVQeQuPXYFdVZ%2F698jFY2qITOJFLZ2JrtQWtUbI5b9%2BWlQu943bx2WHIw8pf%2BwwY6pQEyUXZ3yRsK014j5lyPU44vIh9Es%2BCl2HlMAAgF9rnXZGl4A5QizJJeQmlBOm35r6qUSEnlD%2BpLsIy30l3EZ75cEl0bA%2BT%2B3eUiDI4eZ%2FxbXSrIyQhPNZQyPrj7RRSzPc6L54%2Ff6z4nz4Ce6GKoobKMoCMB7tHaIX62VuTCLlyEfM6GMaTbx5BqwbyGOTj4a4cWcf%2FAJpw3E8Wh3yBiWCharuPNFm0%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2023%20Jul%202025%2000%3A01%3A14%20GMT&X-Amz-Signature=80e52700d6838eb9ecd96ad701af94c858d0fdb2e118e22941c0d8b6fe1b54c2"
$ws.Cells.Item(2, 12).Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991265517/1753110599426/NHvGumSzkA-camera-video-segment-1753110601926.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFFWEYYQK%2F20250722%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250722T160114Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjENX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFA8Y0PgXhDtNu2%2BQHuQpEZwh0JskrXSr5GnqoS%2FtvyLAiEA14WPHxlFKgBdu70KIBqHbLoKFzI6fR80kkWgcedYHasq5gMI7v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDDKCurDpNn5HydYpgCq6A7tVZ4RPeUwBMWtl7mltXrtDVDjafXNo41hLbZ8mNNOiXR0bnL%2FSwLen7qeGKU1jmlnfe0DcMSjLMqX%2FXjtuB%2FBIKxNbO5YvcixVrqwuMVhKdK5s0KbZrzZlSEVpt5ir%2Bb2slBXuGpBE%2Fk1negUCLUvngM1WwVQ9lSLEU5eu75OS6bSJ6ZT%2BmYw9UmY%2BFof9%2FmWjTehN6qLDSMCeo4gcUg2nrH%2Bh620YdZ1sJcTegwmxxf5GU0ZUWuuQRIMRphd5xiv1K%2FJgPiDkAJCAqgBBIs2E5%2FfH6Fu8TXYTl9P%2FI0PT%2FZvy%2BN1i0HMilb6vKDTQsvZ4TcMzmSyVtvSwSaQzzZlpI3KgaEmhUiJKowNpXdgu88rXpBTMhyprbQH8P1MBr2n9ONmLbu0g4RjQZMRfUfey6uAfPbEur77F%2BTjhBEDM3mwejefwA3AeKZv7TiEgvbDOD9QYMxZV8XUTR9kMNlJCsrVR7pDR2CMC9%2F32ALrQIbas6UXYrRPgglwN9uobg59jHNiY72uGCJCe%2Bp0NIF6sVQeQuPXYFdVZ%2F698jFY2qITOJFLZ2JrtQWtUbI5b9%2BWlQu943bx2WHIw8pf%2BwwY6pQEyUXZ3yRsK014j5lyPU44vIh9Es%2BCl2HlMAAgF9rnXZGl4A5QizJJeQmlBOm35r6qUSEnlD%2BpLsIy30l3EZ75cEl0bA%2BT%2B3eUiDI4eZ%2FxbXSrIyQhPNZQyPrj7RRSzPc6L54%2Ff6z4nz4Ce6GKoobKMoCMB7tHaIX62VuTCLlyEfM6GMaTbx5BqwbyGOTj4a4cWcf%2FAJpw3E8Wh3yBiWCharuPNFm0%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2023%20Jul%202025%2000%3A01%3A14%20GMT&X-Amz-Signature=61b8455b473921879f455b8126ce49daa4476845b9627682d3397878a400c2b4"
